# ============================================================================
# Edit: add "2022-Q3" sheet (with fund holdings data) to the workbook,
# positioned right after "总计" and before "2022-Q2"; also update the
# "总计" summary sheet with a new row for 2022-Q3.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3,
#    push 2022-Q2 / 2022-Q1 / 2021-Q4 down by one row.
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 39
$summary.Range("D2").Value = 11.21

# ----------------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet right after "总计".
# ----------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

# Copy header row (values + style) from the existing "2022-Q2" sheet,
# since the header text/formatting is identical across quarterly sheets.
$q2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Copy the index-column style (bold + border, used on A2:A40) from
# the existing "2022-Q2" sheet's A2 cell.
$q2.Range("A2").Copy($newSheet.Range("A2:A40"))

# ----------------------------------------------------------------------
# 3) Populate the fund holdings data rows (row 2 .. row 40).
#    Column A = index (number), B..F = text, G = text (except last row,
#    which is a literal numeric 0), H = number.
# ----------------------------------------------------------------------
$q3rows = @(
    ,@(0, '007130', '中庚小盘价值股票', '75.87', '93.06', '2.89', '2.1926', 8)
    ,@(1, '012930', '中庚价值先锋股票', '47.83', '94.71', '4.49', '2.1476', 7)
    ,@(2, '166301', '华商新趋势优选灵活配置混合', '57.69', '75.38', '1.81', '1.0442', 9)
    ,@(3, '009646', '南方核心成长混合A', '14.92', '85.80', '5.30', '0.7908', 2)
    ,@(4, '000390', '华商优势行业混合', '26.43', '83.59', '2.21', '0.5841', 7)
    ,@(5, '202011', '南方优选价值混合A', '11.30', '83.49', '5.14', '0.5808', 2)
    ,@(6, '010132', '南方创新成长混合A', '8.38', '87.72', '5.34', '0.4475', 2)
    ,@(7, '003378', '泰康策略优选灵活配置混合', '13.78', '83.10', '2.85', '0.3927', 10)
    ,@(8, '009681', '南方创新精选一年定期开放混合A', '6.48', '84.09', '5.87', '0.3804', 1)
    ,@(9, '001150', '融通互联网传媒灵活配置混合', '8.61', '88.88', '4.17', '0.3590', 4)
    ,@(10, '378010', '上投摩根成长先锋混合A', '13.53', '83.14', '2.22', '0.3004', 10)
    ,@(11, '010874', '泰康品质生活混合A', '6.86', '84.42', '2.86', '0.1962', 10)
    ,@(12, '009353', '浙商科技创新一个月滚动持有混合A', '2.06', '91.05', '8.71', '0.1794', 1)
    ,@(13, '166801', '浙商聚潮新思维混合A', '1.89', '78.76', '8.35', '0.1578', 1)
    ,@(14, '009647', '南方核心成长混合C', '2.91', '85.80', '5.30', '0.1542', 2)
    ,@(15, '013091', '上投摩根均衡优选混合A', '6.57', '76.63', '2.17', '0.1426', 9)
    ,@(16, '009682', '南方创新精选一年定期开放混合C', '2.39', '84.09', '5.87', '0.1403', 1)
    ,@(17, '005729', '南方人工智能主题混合', '1.91', '78.93', '6.84', '0.1306', 2)
    ,@(18, '009354', '浙商科技创新一个月滚动持有混合C', '1.45', '91.05', '8.71', '0.1263', 1)
    ,@(19, '002577', '南方新兴龙头灵活配置混合', '1.52', '78.68', '6.69', '0.1017', 1)
    ,@(20, '010875', '泰康品质生活混合C', '3.44', '84.42', '2.86', '0.0984', 10)
    ,@(21, '010133', '南方创新成长混合C', '1.84', '87.72', '5.34', '0.0983', 2)
    ,@(22, '014085', '浙商聚潮新思维混合C', '0.85', '78.76', '8.35', '0.0710', 1)
    ,@(23, '015373', '浙商智选新兴产业混合A', '0.70', '92.12', '8.92', '0.0624', 1)
    ,@(24, '001692', '南方国策动力股票', '2.60', '94.22', '2.01', '0.0523', 10)
    ,@(25, '013899', '上投摩根全景优势股票A', '2.31', '81.01', '2.18', '0.0504', 10)
    ,@(26, '012904', '上投摩根鑫睿优选一年持有期混合', '1.95', '74.92', '2.12', '0.0413', 9)
    ,@(27, '001723', '华商新动力灵活配置混合', '0.70', '84.40', '4.58', '0.0321', 3)
    ,@(28, '012669', '南方新兴产业混合A', '0.53', '81.22', '5.83', '0.0309', 1)
    ,@(29, '002293', '南方益和灵活配置混合', '1.12', '75.22', '2.71', '0.0304', 9)
    ,@(30, '015374', '浙商智选新兴产业混合C', '0.27', '92.12', '8.92', '0.0241', 1)
    ,@(31, '012670', '南方新兴产业混合C', '0.40', '81.22', '5.83', '0.0233', 1)
    ,@(32, '013092', '上投摩根均衡优选混合C', '0.63', '76.63', '2.17', '0.0137', 9)
    ,@(33, '013903', '国泰君安信息行业混合', '0.25', '71.75', '5.04', '0.0126', 2)
    ,@(34, '006539', '南方优选价值混合C', '0.15', '83.49', '5.14', '0.0077', 2)
    ,@(35, '002123', '北信瑞丰外延增长主题灵活配置混合', '0.15', '88.55', '4.37', '0.0066', 6)
    ,@(36, '960020', '南方优选价值混合H', '0.08', '83.49', '5.14', '0.0041', 2)
    ,@(37, '013900', '上投摩根全景优势股票C', '0.14', '81.01', '2.18', '0.0031', 10)
    ,@(38, '015077', '上投摩根成长先锋混合C', '0.00', '83.14', '2.22', '0', 10)
)

# Columns B:G hold values that look numeric (e.g. "75.87") but must be
# stored as text, matching the source data. Force text interpretation by
# setting the number format to Text before assigning the values.
$dataRange = $newSheet.Range("B2:G40")
$dataRange.NumberFormat = "@"

$r = 2
foreach ($row in $q3rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the Text number-format we applied above (keeps cells plain/"General"
# like the rest of the workbook) without disturbing the text values already
# stored in them.
$dataRange.ClearFormats()
$q2.Range("A2").Copy($newSheet.Range("A2:A40"))

# Row 40's "持有市值(亿元)" (column G) is a literal numeric 0 in the source
# data (not the text "0.0000"), so restore it to a real number.
$newSheet.Range("G40").Value = 0

Write-Host "2022-Q3 sheet added and 总计 updated"
